# Auto-generated edit script
# Applies numeric value updates (and a few cell additions/removals)
# to the 'Typhon_Profits' market-data workbook, per scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1066.6666
$ws.Range("J32").Value = 1066.6666
$ws.Range("L32").Value = 1066.6666
$ws.Range("N32").Value = -1718.6666
$ws.Range("H40").Value = 1591.6666
$ws.Range("I40").Value = 940
$ws.Range("J40").Value = 2057.1428
$ws.Range("K40").Value = 940
$ws.Range("L40").Value = 2057.1428
$ws.Range("M40").Value = -765
$ws.Range("N40").Value = -2407.1428
$ws.Range("H106").Value = 1941.04
$ws.Range("I106").Value = 1150.2858
$ws.Range("J106").Value = 2947.4546
$ws.Range("K106").Value = 1150.2858
$ws.Range("L106").Value = 2947.4546
$ws.Range("M106").Value = -519.2858000000001
$ws.Range("N106").Value = -4209.4546
$ws.Range("H132").Value = 26438.977
$ws.Range("I132").Value = 26438.977
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 79316.931
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -76786.931
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 2675.0557
$ws.Range("I138").Value = 1660.2941
$ws.Range("J138").Value = 3583
$ws.Range("K138").Value = 4980.8823
$ws.Range("L138").Value = 10749
$ws.Range("M138").Value = 159.1176999999998
$ws.Range("N138").Value = -21029
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2221.8386
$ws.Range("I61").Value = 1121.3529
$ws.Range("J61").Value = 3558.1428
$ws.Range("K61").Value = 1121.3529
$ws.Range("L61").Value = 3558.1428
$ws.Range("M61").Value = -909.3529000000001
$ws.Range("N61").Value = -3982.1428
$ws.Range("H63").Value = 2345.111
$ws.Range("I63").Value = 2362.5
$ws.Range("K63").Value = 2362.5
$ws.Range("M63").Value = -1676.5
$ws.Range("H66").Value = 2345.111
$ws.Range("I66").Value = 2362.5
$ws.Range("K66").Value = 11812.5
$ws.Range("M66").Value = -8380.5
$ws.Range("H132").Value = 22689.041
$ws.Range("I132").Value = 1484.7222
$ws.Range("K132").Value = 4454.1666
$ws.Range("M132").Value = -1924.1666
$ws.Range("H135").Value = 26318.428
$ws.Range("J135").Value = 26318.428
$ws.Range("L135").Value = 26318.428
$ws.Range("N135").Value = -36458.428
$ws.Range("H136").Value = 2221.8386
$ws.Range("I136").Value = 1121.3529
$ws.Range("J136").Value = 3558.1428
$ws.Range("K136").Value = 3364.0587
$ws.Range("L136").Value = 10674.4284
$ws.Range("M136").Value = -814.0587000000005
$ws.Range("N136").Value = -15774.4284
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2961.0667
$ws.Range("I94").Value = 1278.4445
$ws.Range("K94").Value = 1278.4445
$ws.Range("M94").Value = -827.4445000000001
$ws.Range("H99").Value = 1749.875
$ws.Range("I99").Value = 1571.4286
$ws.Range("K99").Value = 1571.4286
$ws.Range("M99").Value = -73.42859999999996
$ws.Range("H107").Value = 1671
$ws.Range("I107").Value = 1716.1666
$ws.Range("K107").Value = 1716.1666
$ws.Range("M107").Value = 203.8334
$ws.Range("H119").Value = 19253.334
$ws.Range("J119").Value = 19253.334
$ws.Range("L119").Value = 19253.334
$ws.Range("N119").Value = -28929.334
$ws.Range("H134").Value = 3493.6428
$ws.Range("I134").Value = 3639.76
$ws.Range("J134").Value = 2276
$ws.Range("K134").Value = 10919.28
$ws.Range("L134").Value = 6828
$ws.Range("M134").Value = -8384.280000000001
$ws.Range("N134").Value = -11898
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 462.33334
$ws.Range("I22").Value = 422.54544
$ws.Range("K22").Value = 422.54544
$ws.Range("M22").Value = -72.54543999999999
$ws.Range("H31").Value = 13108.6
$ws.Range("I31").Value = 26226.25
$ws.Range("K31").Value = 26226.25
$ws.Range("M31").Value = -25931.25
$ws.Range("H34").Value = 13108.6
$ws.Range("I34").Value = 26226.25
$ws.Range("K34").Value = 26226.25
$ws.Range("M34").Value = -26024.25
$ws.Range("H58").Value = 11808.218
$ws.Range("J58").Value = 56288.89
$ws.Range("L58").Value = 56288.89
$ws.Range("N58").Value = -56694.89
$ws.Range("H99").Value = 5604.737
$ws.Range("I99").Value = 4179
$ws.Range("J99").Value = 7188.8887
$ws.Range("K99").Value = 4179
$ws.Range("L99").Value = 7188.8887
$ws.Range("M99").Value = -2681
$ws.Range("N99").Value = -10184.8887
$ws.Range("H126").Value = 5604.737
$ws.Range("I126").Value = 4179
$ws.Range("J126").Value = 7188.8887
$ws.Range("K126").Value = 12537
$ws.Range("L126").Value = 21566.6661
$ws.Range("M126").Value = -10067
$ws.Range("N126").Value = -26506.6661
$ws.Range("H132").Value = 23025.92
$ws.Range("J132").Value = 4381.0713
$ws.Range("L132").Value = 13143.2139
$ws.Range("N132").Value = -18203.2139
$ws.Range("H134").Value = 1214.2858
$ws.Range("I134").Value = 1250
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 3750
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -1215
$ws.Range("N134").Value = -8070
$ws.Range("H136").Value = 11808.218
$ws.Range("J136").Value = 56288.89
$ws.Range("L136").Value = 168866.67
$ws.Range("N136").Value = -173966.67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3835.5386
$ws.Range("J68").Value = 5285.852
$ws.Range("L68").Value = 15857.556
$ws.Range("N68").Value = -17479.556
$ws.Range("H71").Value = 3835.5386
$ws.Range("J71").Value = 5285.852
$ws.Range("L71").Value = 47572.668
$ws.Range("N71").Value = -55684.668
$ws.Range("H87").Value = 19250
$ws.Range("I87").Value = 9666.666999999999
$ws.Range("K87").Value = 29000.001
$ws.Range("M87").Value = -27752.001
$ws.Range("H90").Value = 19250
$ws.Range("I90").Value = 9666.666999999999
$ws.Range("K90").Value = 87000.003
$ws.Range("M90").Value = -80760.003
$ws.Range("H98").Value = 750
$ws.Range("J98").Value = 750
$ws.Range("L98").Value = 2250
$ws.Range("N98").Value = -5246
$ws.Range("H129").Value = 1019.4286
$ws.Range("I129").Value = 427.2
$ws.Range("K129").Value = 1281.6
$ws.Range("M129").Value = 3718.4
$ws.Range("H131").Value = 786.48486
$ws.Range("J131").Value = 803.6667
$ws.Range("L131").Value = 2411.0001
$ws.Range("N131").Value = -12491.0001
$ws.Range("H140").Value = 5731.2593
$ws.Range("I140").Value = 7521.75
$ws.Range("J140").Value = 3126.9092
$ws.Range("K140").Value = 22565.25
$ws.Range("L140").Value = 9380.7276
$ws.Range("M140").Value = -17385.25
$ws.Range("N140").Value = -19740.7276
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 27663
$ws.Range("J125").Value = 27663
$ws.Range("L125").Value = 27663
$ws.Range("N125").Value = -32583
$ws.Range("H132").Value = 17995.088
$ws.Range("I132").Value = 3476.4583
$ws.Range("J132").Value = 52839.8
$ws.Range("K132").Value = 10429.3749
$ws.Range("L132").Value = 158519.4
$ws.Range("M132").Value = -7899.374899999999
$ws.Range("N132").Value = -163579.4
$ws.Range("H136").Value = 14175.333
$ws.Range("J136").Value = 14175.333
$ws.Range("L136").Value = 42525.999
$ws.Range("N136").Value = -47625.999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4300.3335
$ws.Range("J22").Value = 2500
$ws.Range("L22").Value = 2500
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 4300.3335
$ws.Range("J27").Value = 2500
$ws.Range("L27").Value = 2500
$ws.Range("N27").Value = -2714
$ws.Range("H46").Value = 1833.1177
$ws.Range("I46").Value = 1690.8667
$ws.Range("K46").Value = 1690.8667
$ws.Range("M46").Value = -1502.8667
$ws.Range("H122").Value = 3513.8823
$ws.Range("I122").Value = 3002.6
$ws.Range("J122").Value = 4244.2856
$ws.Range("K122").Value = 9007.799999999999
$ws.Range("L122").Value = 12732.8568
$ws.Range("M122").Value = -6557.799999999999
$ws.Range("N122").Value = -17632.8568
$ws.Range("H132").Value = 1600.2941
$ws.Range("I132").Value = 1058.8422
$ws.Range("J132").Value = 2286.1333
$ws.Range("K132").Value = 3176.5266
$ws.Range("L132").Value = 6858.3999
$ws.Range("M132").Value = -646.5266000000001
$ws.Range("N132").Value = -11918.3999
$ws.Range("H136").Value = 2334.8708
$ws.Range("I136").Value = 1193.5264
$ws.Range("J136").Value = 4142
$ws.Range("K136").Value = 3580.5792
$ws.Range("L136").Value = 12426
$ws.Range("M136").Value = -1030.5792
$ws.Range("N136").Value = -17526
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H126").Value = 1386.8
$ws.Range("J126").Value = 1992.75
$ws.Range("L126").Value = 5978.25
$ws.Range("N126").Value = -10918.25
$ws.Range("H137").Value = 41536.25
$ws.Range("J137").Value = 41536.25
$ws.Range("L137").Value = 41536.25
$ws.Range("N137").Value = -51736.25
